# Applies the Sun Sep 3 2023 cryptos-list data refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping it stored as plain text
# (the sheet stores prices like "1.003" / "215.06" as text, not numbers).
function Set-TextValue([string]$addr, [string]$text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '25.934.64'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.639.68'
$ws.Range("E3").Value = '  +0.06%  '
Set-TextValue "D4" '1.003'
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue "D5" '215.06'
$ws.Range("E5").Value = '  -0.01%  '
Set-TextValue "D6" '0.5059'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.20%  '
Set-TextValue "D8" '0.2562'
$ws.Range("E8").Value = '  -0.52%  '
Set-TextValue "D9" '0.06368'
$ws.Range("E9").Value = '  -0.12%  '
Set-TextValue "D10" '19.47'
$ws.Range("E10").Value = '  -0.29%  '
Set-TextValue "D11" '0.07745'
$ws.Range("E11").Value = '  -0.01%  '
Set-TextValue "D12" '4.283'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = '1.646.43'
$ws.Range("E13").Value = '  +0.27%  '
Set-TextValue "D14" '0.5445'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '0.0₅7827'
$ws.Range("E15").Value = '  -0.87%  '
Set-TextValue "D16" '64.21'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '25.975.10'
$ws.Range("E17").Value = '  +0.35%  '
Set-TextValue "D18" '1.003'
$ws.Range("E18").Value = '  +0.07%  '
Set-TextValue "D19" '197.35'
$ws.Range("E19").Value = '  -2.62%  '
Set-TextValue "D20" '4.438'
$ws.Range("E20").Value = '  +0.96%  '
Set-TextValue "D21" '9.941'
$ws.Range("E21").Value = '  +0.69%  '
Set-TextValue "D22" '6.046'
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("E23").Value = '  +0.28%  '
Set-TextValue "D24" '1.899'
$ws.Range("E24").Value = '  +1.65%  '
Set-TextValue "D25" '140.99'
Set-TextValue "D26" '0.1168'
$ws.Range("E26").Value = '  +3.10%  '
Set-TextValue "D27" '6.871'
$ws.Range("E27").Value = '  +1.33%  '
Set-TextValue "D28" '15.69'
$ws.Range("E28").Value = '  +0.11%  '
Set-TextValue "D29" '1.236'
$ws.Range("E29").Value = '  -0.44%  '
Set-TextValue "D30" '0.04976'
$ws.Range("E30").Value = '  +0.09%  '
Set-TextValue "D31" '3.258'
$ws.Range("E31").Value = '  -0.45%  '
Set-TextValue "D32" '3.181'
$ws.Range("E32").Value = '  -0.39%  '
Set-TextValue "D33" '1.539'
$ws.Range("E33").Value = '  -0.20%  '
Set-TextValue "D34" '2.360'
$ws.Range("E34").Value = '  -0.23%  '
Set-TextValue "D35" '0.8926'
$ws.Range("E35").Value = '  +0.28%  '
Set-TextValue "D36" '2.594'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").Value = '1.131.36'
$ws.Range("E37").Value = '  -1.82%  '
Set-TextValue "D38" '0.5439'
$ws.Range("E38").Value = '  -2.94%  '
Set-TextValue "D39" '0.01554'
$ws.Range("E39").Value = '  -0.56%  '
Set-TextValue "D40" '2.553'
$ws.Range("E40").Value = '  -0.30%  '
Set-TextValue "D41" '1.003'
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").Value = '0.0₈129'
$ws.Range("E42").Value = '  +10.73%  '
Set-TextValue "D43" '5.594'
$ws.Range("E43").Value = '  -1.53%  '
Set-TextValue "D44" '0.8140'
Set-TextValue "D45" '99.64'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '1.778.31'
$ws.Range("E46").Value = '  +0.12%  '
Set-TextValue "D47" '0.4538'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("E48").Value = '  -0.32%  '
Set-TextValue "D49" '54.71'
$ws.Range("E49").Value = '  -0.11%  '
Set-TextValue "D50" '0.05070'
$ws.Range("E50").Value = '  +0.35%  '
Set-TextValue "D51" '1.005'
$ws.Range("E51").Value = '  +0.59%  '
